$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 17.77319415614589
$ws.Range("C2").Value = 7.791983378675867
$ws.Range("D2").Value = 8.239479812129757
$ws.Range("E2").Value = 12.36374021842432
$ws.Range("F2").Value = 33.83264864516143
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 24.73480789260664
$ws.Range("J2").Value = 9.581487980166377
$ws.Range("L2").Value = 11.48276104519301
$ws.Range("N2").Value = 18.03796723850543
$ws.Range("O2").Value = 25.90925122623008
# Row 3
$ws.Range("B3").Value = 17.34527673506555
$ws.Range("C3").Value = 7.480115152201347
$ws.Range("D3").Value = 8.239041773105118
$ws.Range("E3").Value = 12.39110689749562
$ws.Range("F3").Value = 33.886380937281
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 24.83563599051525
$ws.Range("J3").Value = 9.601598045439776
$ws.Range("L3").Value = 11.46747575642975
$ws.Range("N3").Value = 18.08465470113766
$ws.Range("O3").Value = 25.97065680667485
# Row 4
$ws.Range("B4").Value = 17.07936528284254
$ws.Range("C4").Value = 7.280451675780484
$ws.Range("D4").Value = 8.2396365845633
$ws.Range("E4").Value = 12.40921478516419
$ws.Range("F4").Value = 33.92738678865064
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 24.90240761032862
$ws.Range("J4").Value = 9.614612442509245
$ws.Range("L4").Value = 11.45953316998241
$ws.Range("N4").Value = 18.11510521220666
$ws.Range("O4").Value = 26.01416566149595
# Row 5
$ws.Range("B5").Value = 16.97036124088529
$ws.Range("C5").Value = 7.197105834763626
$ws.Range("D5").Value = 8.240097022750712
$ws.Range("E5").Value = 12.4169224431989
$ws.Range("F5").Value = 33.94610871400028
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 24.93083900411917
$ws.Range("J5").Value = 9.620084030234262
$ws.Range("L5").Value = 11.45666184680705
$ws.Range("N5").Value = 18.12796355970592
$ws.Range("O5").Value = 26.03335261239136
# Row 6
$ws.Range("B6").Value = 16.95222699757618
$ws.Range("C6").Value = 7.183149023794773
$ws.Range("D6").Value = 8.240186671439924
$ws.Range("E6").Value = 12.41822215161721
$ws.Range("F6").Value = 33.94933886299014
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 24.93563374640776
$ws.Range("J6").Value = 9.621002750475997
$ws.Range("L6").Value = 11.45620720349888
$ws.Range("N6").Value = 18.13012585379877
$ws.Range("O6").Value = 26.03662647772526
# Row 7
$ws.Range("B7").Value = 17.07789761132723
$ws.Range("C7").Value = 7.279335561414098
$ws.Range("D7").Value = 8.239641910239062
$ws.Range("E7").Value = 12.40931740241263
$ws.Range("F7").Value = 33.92763113886468
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 24.90278610238473
$ws.Range("J7").Value = 9.614685552951316
$ws.Range("L7").Value = 11.45949296377202
$ws.Range("N7").Value = 18.11527680315364
$ws.Range("O7").Value = 26.01441852949623
# Row 8
$ws.Range("B8").Value = 17.6263958784398
$ws.Range("C8").Value = 7.686189513204285
$ws.Range("D8").Value = 8.239150021738109
$ws.Range("E8").Value = 12.37290576458991
$ws.Range("F8").Value = 33.84951033247238
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 24.76856333501053
$ws.Range("J8").Value = 9.588283826600152
$ws.Range("L8").Value = 11.47719285776904
$ws.Range("N8").Value = 18.05369527576083
$ws.Range("O8").Value = 25.92921712678445
# Row 9
$ws.Range("B9").Value = 18.67048084795387
$ws.Range("C9").Value = 8.416382288640852
$ws.Range("D9").Value = 8.245000668955846
$ws.Range("E9").Value = 12.3118326971027
$ws.Range("F9").Value = 33.76003277748742
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 24.5439976811214
$ws.Range("J9").Value = 9.541779362619829
$ws.Range("L9").Value = 11.52322913442226
$ws.Range("N9").Value = 17.94705256305001
$ws.Range("O9").Value = 25.80833122873828
# Row 10
$ws.Range("B10").Value = 19.41039411676416
$ws.Range("C10").Value = 8.908498719993192
$ws.Range("D10").Value = 8.253400555666472
$ws.Range("E10").Value = 12.27323016766596
$ws.Range("F10").Value = 33.73327545391898
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 24.40265162616973
$ws.Range("J10").Value = 9.510795320148734
$ws.Range("L10").Value = 11.56378698384855
$ws.Range("N10").Value = 17.8772563251918
$ws.Range("O10").Value = 25.74783805944563
# Row 11
$ws.Range("B11").Value = 19.73959102735281
$ws.Range("C11").Value = 9.122235878417872
$ws.Range("D11").Value = 8.258100060962262
$ws.Range("E11").Value = 12.25702360256403
$ws.Range("F11").Value = 33.72958175584904
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 24.34350166672222
$ws.Range("J11").Value = 9.497384728146557
$ws.Range("L11").Value = 11.58366025738696
$ws.Range("N11").Value = 17.8473507555464
$ws.Range("O11").Value = 25.72649582001914
# Row 12
$ws.Range("B12").Value = 19.86307287512811
$ws.Range("C12").Value = 9.201681942856727
$ws.Range("D12").Value = 8.260004786324702
$ws.Range("E12").Value = 12.25108080249982
$ws.Range("F12").Value = 33.72940198988085
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 24.32184518108787
$ws.Range("J12").Value = 9.492404404612618
$ws.Range("L12").Value = 11.59138670365392
$ws.Range("N12").Value = 17.8362908252081
$ws.Range("O12").Value = 25.71930390462457
# Row 13
$ws.Range("B13").Value = 19.83653296497061
$ws.Range("C13").Value = 9.184638648298359
$ws.Range("D13").Value = 8.259589024100281
$ws.Range("E13").Value = 12.25235205706768
$ws.Range("F13").Value = 33.72938650449841
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 24.32647624764397
$ws.Range("J13").Value = 9.493472655362998
$ws.Range("L13").Value = 11.58971380390688
$ws.Range("N13").Value = 17.8386610185955
$ws.Range("O13").Value = 25.72081320837629
# Row 14
$ws.Range("B14").Value = 19.74977418406463
$ws.Range("C14").Value = 9.128802028067559
$ws.Range("D14").Value = 8.258254263146807
$ws.Range("E14").Value = 12.25653079366867
$ws.Range("F14").Value = 33.72954254085732
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 24.34170508654625
$ws.Range("J14").Value = 9.496973032675475
$ws.Range("L14").Value = 11.58429191345824
$ws.Range("N14").Value = 17.84643554796172
$ws.Range("O14").Value = 25.72588629255145
# Row 15
$ws.Range("B15").Value = 19.69647528380892
$ws.Range("C15").Value = 9.094405295884549
$ws.Range("D15").Value = 8.257452942421763
$ws.Range("E15").Value = 12.25911567836882
$ws.Range("F15").Value = 33.72979684130931
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 24.35112991839605
$ws.Range("J15").Value = 9.499129862164917
$ws.Range("L15").Value = 11.58099689224007
$ws.Range("N15").Value = 17.85123211891641
$ws.Range("O15").Value = 25.72910964396941
# Row 16
$ws.Range("B16").Value = 19.38872153362163
$ws.Range("C16").Value = 8.894323550732866
$ws.Range("D16").Value = 8.253111006965923
$ws.Range("E16").Value = 12.27431651343293
$ws.Range("F16").Value = 33.73368745186004
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 24.40662098891087
$ws.Range("J16").Value = 9.511685466576511
$ws.Range("L16").Value = 11.56251653840194
$ws.Range("N16").Value = 17.87924779537816
$ws.Range("O16").Value = 25.74935726972836
# Row 17
$ws.Range("B17").Value = 19.19794640109052
$ws.Range("C17").Value = 8.768958874030867
$ws.Range("D17").Value = 8.250671515118064
$ws.Range("E17").Value = 12.28398820746204
$ws.Range("F17").Value = 33.73824582807361
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 24.44198310502356
$ws.Range("J17").Value = 9.5195628705185
$ws.Range("L17").Value = 11.5515413351168
$ws.Range("N17").Value = 17.89690659975565
$ws.Range("O17").Value = 25.76336183433694
# Row 18
$ws.Range("B18").Value = 19.08752796346528
$ws.Range("C18").Value = 8.695900743450233
$ws.Range("D18").Value = 8.24935113313332
$ws.Range("E18").Value = 12.28967857230646
$ws.Range("F18").Value = 33.7416657720417
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 24.4628069998819
$ws.Range("J18").Value = 9.524158169785807
$ws.Range("L18").Value = 11.54536289552818
$ws.Range("N18").Value = 17.90723718373091
$ws.Range("O18").Value = 25.7719981596914
# Row 19
$ws.Range("B19").Value = 19.05002738069459
$ws.Range("C19").Value = 8.671002221977362
$ws.Range("D19").Value = 8.248918321807814
$ws.Range("E19").Value = 12.2916271365527
$ws.Range("F19").Value = 33.7429607716461
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 24.46994077586892
$ws.Range("J19").Value = 9.525725136788534
$ws.Range("L19").Value = 11.5432941451058
$ws.Range("N19").Value = 17.91076479676973
$ws.Range("O19").Value = 25.77502204483387
# Row 20
$ws.Range("B20").Value = 19.21832697375977
$ws.Range("C20").Value = 8.782402924843295
$ws.Range("D20").Value = 8.250922647344678
$ws.Range("E20").Value = 12.28294545028619
$ws.Range("F20").Value = 33.73767798112767
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 24.4381685855485
$ws.Range("J20").Value = 9.518717642679464
$ws.Range("L20").Value = 11.55269580128222
$ws.Range("N20").Value = 17.89500881673666
$ws.Range("O20").Value = 25.76181085122947
# Row 21
$ws.Range("B21").Value = 19.77529017692301
$ws.Range("C21").Value = 9.145243329244696
$ws.Range("D21").Value = 8.258642928642601
$ws.Range("E21").Value = 12.25529812821113
$ws.Range("F21").Value = 33.72946363287272
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 24.33721184592988
$ws.Range("J21").Value = 9.495942231207833
$ws.Range("L21").Value = 11.58587903375432
$ws.Range("N21").Value = 17.84414480260944
$ws.Range("O21").Value = 25.72437204015516
# Row 22
$ws.Range("B22").Value = 20.13238550825284
$ws.Range("C22").Value = 9.373674055645363
$ws.Range("D22").Value = 8.26441731342835
$ws.Range("E22").Value = 12.23836119119736
$ws.Range("F22").Value = 33.73119941134077
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 24.27555847550051
$ws.Range("J22").Value = 9.481628077973124
$ws.Range("L22").Value = 11.60873512541824
$ws.Range("N22").Value = 17.81244457652715
$ws.Range("O22").Value = 25.7050913435649
# Row 23
$ws.Range("B23").Value = 19.94246501385172
$ws.Range("C23").Value = 9.25256309808098
$ws.Range("D23").Value = 8.261269141226578
$ws.Range("E23").Value = 12.2472972960573
$ws.Range("F23").Value = 33.72962323289711
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 24.30806742801812
$ws.Range("J23").Value = 9.489215708822121
$ws.Range("L23").Value = 11.5964307428929
$ws.Range("N23").Value = 17.82922267225662
$ws.Range("O23").Value = 25.71490665020633
# Row 24
$ws.Range("B24").Value = 19.20911521068524
$ws.Range("C24").Value = 8.776327932704849
$ws.Range("D24").Value = 8.250808854610799
$ws.Range("E24").Value = 12.28341647590787
$ws.Range("F24").Value = 33.73793221502793
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 24.43989159194381
$ws.Range("J24").Value = 9.519099563090199
$ws.Range("L24").Value = 11.55217345770121
$ws.Range("N24").Value = 17.89586624897869
$ws.Range("O24").Value = 25.76251022883798
# Row 25
$ws.Range("B25").Value = 18.3922583709789
$ws.Range("C25").Value = 8.22644019292348
$ws.Range("D25").Value = 8.242693582126648
$ws.Range("E25").Value = 12.32725174165337
$ws.Range("F25").Value = 33.77740024798904
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 24.60060289914222
$ws.Range("J25").Value = 9.55379900650971
$ws.Range("L25").Value = 11.50957929675195
$ws.Range("N25").Value = 17.97439630124532
$ws.Range("O25").Value = 25.83607082150786
